$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume/Change% (E) columns per latest data refresh.
# D-column values are forced to remain text (matching the sheet's existing inlineStr cells)
# by prefixing with an apostrophe and resetting the cell style afterwards, since some of the
# prices (e.g. "150.95", "1.00") would otherwise be auto-parsed by Excel as numbers.

$cell = $ws.Range("D2")
$cell.Value = "'63.567.54"
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +2.34%  '
$cell = $ws.Range("D3")
$cell.Value = "'2.565.13"
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +5.22%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  +3.00%  '
$cell = $ws.Range("D6")
$cell.Value = "'150.95"
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +8.91%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +1.00%  '
$cell = $ws.Range("D9")
$cell.Value = "'2.563.82"
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +5.30%  '
$ws.Range("E10").Value = '  +2.38%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("E13").Value = '  +3.56%  '
$cell = $ws.Range("D14")
$cell.Value = "'28.22"
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +9.11%  '
$cell = $ws.Range("D15")
$cell.Value = "'3.026.96"
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +5.38%  '
$cell = $ws.Range("D16")
$cell.Value = "'63.543.33"
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +2.44%  '
$ws.Range("E17").Value = '  +2.82%  '
$cell = $ws.Range("D18")
$cell.Value = "'2.561.53"
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +5.12%  '
$cell = $ws.Range("D20")
$cell.Value = "'343.02"
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.80%  '
$cell = $ws.Range("D21")
$cell.Value = "'4.40"
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +3.83%  '
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("E23").Value = '  +0.25%  '
$cell = $ws.Range("D24")
$cell.Value = "'66.26"
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("E25").Value = '  -1.28%  '
$cell = $ws.Range("D26")
$cell.Value = "'1.61"
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +5.95%  '
$cell = $ws.Range("D27")
$cell.Value = "'8.55"
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +3.25%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  +8.77%  '
$cell = $ws.Range("D30")
$cell.Value = "'7.15"
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +12.55%  '
$cell = $ws.Range("D31")
$cell.Value = "'0.0₃0844"
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +7.06%  '
$ws.Range("E32").Value = '  +3.78%  '
$cell = $ws.Range("D33")
$cell.Value = "'177.41"
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +3.64%  '
$cell = $ws.Range("D34")
$cell.Value = "'1.57"
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +8.10%  '
$cell = $ws.Range("D35")
$cell.Value = "'421.24"
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +12.37%  '
$cell = $ws.Range("D36")
$cell.Value = "'0.409"
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +3.03%  '
$cell = $ws.Range("D37")
$cell.Value = "'19.22"
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +3.44%  '
$cell = $ws.Range("D38")
$cell.Value = "'4.47"
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -0.14%  '
$cell = $ws.Range("D40")
$cell.Value = "'1.75"
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +3.75%  '
$cell = $ws.Range("D41")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +0.04%  '
$cell = $ws.Range("D42")
$cell.Value = "'40.15"
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +2.10%  '
$cell = $ws.Range("D43")
$cell.Value = "'156.30"
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +6.57%  '
$cell = $ws.Range("D44")
$cell.Value = "'3.83"
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +4.36%  '
$cell = $ws.Range("D45")
$cell.Value = "'21.33"
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +2.72%  '
$cell = $ws.Range("D46")
$cell.Value = "'0.613"
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +4.01%  '
$ws.Range("E47").Value = '  +3.13%  '
$ws.Range("E48").Value = '  +1.11%  '
$cell = $ws.Range("D49")
$cell.Value = "'0.0236"
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +6.15%  '
$cell = $ws.Range("D50")
$cell.Value = "'18.92"
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +5.09%  '
$ws.Range("E51").Value = '  +9.02%  '
